$wb = $excel.ActiveWorkbook

# Hyperlink colour used throughout this workbook for the "file name" style
# links (matches existing FF6495ED, i.e. RGB(0x64,0x95,0xED) in COM BGR packing).
$linkColor = 15570276

function Set-LinkFormat($rng) {
    $rng.Font.Underline = 2
    $rng.Font.Color = $linkColor
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Drop existing hyperlinks on the sheet so we can rebuild them cleanly.
$ws.Range("A1").Hyperlinks.Delete()

# The row that used to sit at A3 (".localization-config") moves to A4.
$ws.Range("A4").Value2 = ".localization-config"
$ws.Range("B4").Value2 = "Not to be localized"
$ws.Range("C4").Value2 = "Not to be localized"

# The previously-handed-off file is now ready again under a new name.
$ws.Range("A2").Value2 = "fe219963-43f6-4357-a768-730c515b9e56.md"
$ws.Range("B2").Value2 = "Ready for handoff"
$ws.Range("C2").Value2 = "Ready for handoff"

# A brand new file shows up as row 3.
$ws.Range("A3").Value2 = "ffff2d99924b-233c-4d18-be84-ddfba0888239.md"
$ws.Range("B3").Value2 = "Ready for handoff"
$ws.Range("C3").Value2 = "Ready for handoff"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e6dd9c9c4069ad55807f26b3fec2c34f20c368b8/e2e/fe219963-43f6-4357-a768-730c515b9e56.md", "", "", "fe219963-43f6-4357-a768-730c515b9e56.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e6dd9c9c4069ad55807f26b3fec2c34f20c368b8/e2e/ffff2d99924b-233c-4d18-be84-ddfba0888239.md", "", "", "ffff2d99924b-233c-4d18-be84-ddfba0888239.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4dc906a52e24d7ef983af6a0794887a1e5a68252/.localization-config", "", "", ".localization-config")

Set-LinkFormat($ws.Range("A2"))
Set-LinkFormat($ws.Range("A3"))
Set-LinkFormat($ws.Range("A4"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A1").Hyperlinks.Delete()

# Old row 3 (".localization-config", "Not to be localized") -> row 4.
$ws.Range("A4").Value2 = ".localization-config"
$ws.Range("B4").Value2 = "Not to be localized"
$ws.Range("D4").Value2 = "0001-01-01 00:00:00"
$ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G4").Value2 = "0001-01-01 00:00:00"
$ws.Range("H4").Value2 = "Ignored"

# Row 2: same source file, now ready, with real handoff info.
$ws.Range("A2").Value2 = "fe219963-43f6-4357-a768-730c515b9e56.md"
$ws.Range("B2").Value2 = "Ready for handoff"
$ws.Range("C2").Value2 = "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf"
$ws.Range("D2").Value2 = "2016-01-20 03:51:39"
$ws.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G2").Value2 = "0001-01-01 00:00:00"
$ws.Range("H2").Value2 = "Include"

# Row 3: brand new source file with the same handoff package.
$ws.Range("A3").Value2 = "ffff2d99924b-233c-4d18-be84-ddfba0888239.md"
$ws.Range("B3").Value2 = "Ready for handoff"
$ws.Range("C3").Value2 = "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf"
$ws.Range("D3").Value2 = "2016-01-20 03:51:39"
$ws.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G3").Value2 = "0001-01-01 00:00:00"
$ws.Range("H3").Value2 = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e6dd9c9c4069ad55807f26b3fec2c34f20c368b8/e2e/fe219963-43f6-4357-a768-730c515b9e56.md", "", "", "fe219963-43f6-4357-a768-730c515b9e56.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTest/oltest/blob/465c21a7bc3af6829ddc588a5812abf7541dc467/e2e/loc/zh-cn/fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf", "", "", "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e6dd9c9c4069ad55807f26b3fec2c34f20c368b8/e2e/ffff2d99924b-233c-4d18-be84-ddfba0888239.md", "", "", "ffff2d99924b-233c-4d18-be84-ddfba0888239.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTest/oltest/blob/465c21a7bc3af6829ddc588a5812abf7541dc467/e2e/loc/zh-cn/fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf", "", "", "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4dc906a52e24d7ef983af6a0794887a1e5a68252/.localization-config", "", "", ".localization-config")

Set-LinkFormat($ws.Range("A2"))
Set-LinkFormat($ws.Range("C2"))
Set-LinkFormat($ws.Range("A3"))
Set-LinkFormat($ws.Range("C3"))
Set-LinkFormat($ws.Range("A4"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A1").Hyperlinks.Delete()

$ws.Range("A4").Value2 = ".localization-config"
$ws.Range("B4").Value2 = "Not to be localized"
$ws.Range("D4").Value2 = "0001-01-01 00:00:00"
$ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G4").Value2 = "0001-01-01 00:00:00"
$ws.Range("H4").Value2 = "Ignored"

$ws.Range("A2").Value2 = "fe219963-43f6-4357-a768-730c515b9e56.md"
$ws.Range("B2").Value2 = "Ready for handoff"
$ws.Range("C2").Value2 = "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf"
$ws.Range("D2").Value2 = "2016-01-20 03:51:49"
$ws.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G2").Value2 = "0001-01-01 00:00:00"
$ws.Range("H2").Value2 = "Include"

$ws.Range("A3").Value2 = "ffff2d99924b-233c-4d18-be84-ddfba0888239.md"
$ws.Range("B3").Value2 = "Ready for handoff"
$ws.Range("C3").Value2 = "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf"
$ws.Range("D3").Value2 = "2016-01-20 03:51:49"
$ws.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G3").Value2 = "0001-01-01 00:00:00"
$ws.Range("H3").Value2 = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e6dd9c9c4069ad55807f26b3fec2c34f20c368b8/e2e/fe219963-43f6-4357-a768-730c515b9e56.md", "", "", "fe219963-43f6-4357-a768-730c515b9e56.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTest/oltest/blob/465c21a7bc3af6829ddc588a5812abf7541dc467/e2e/loc/de-de/fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf", "", "", "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e6dd9c9c4069ad55807f26b3fec2c34f20c368b8/e2e/ffff2d99924b-233c-4d18-be84-ddfba0888239.md", "", "", "ffff2d99924b-233c-4d18-be84-ddfba0888239.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTest/oltest/blob/465c21a7bc3af6829ddc588a5812abf7541dc467/e2e/loc/de-de/fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf", "", "", "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4dc906a52e24d7ef983af6a0794887a1e5a68252/.localization-config", "", "", ".localization-config")

Set-LinkFormat($ws.Range("A2"))
Set-LinkFormat($ws.Range("C2"))
Set-LinkFormat($ws.Range("A3"))
Set-LinkFormat($ws.Range("C3"))
Set-LinkFormat($ws.Range("A4"))
